$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Cells.Item(2,4) '68.763.57'
$ws.Cells.Item(2,5).Value2 = '  -0.44%  '

Set-TextValue $ws.Cells.Item(3,4) '2.429.51'
$ws.Cells.Item(3,5).Value2 = '  -1.91%  '

$ws.Cells.Item(4,5).Value2 = '  +0.02%  '

Set-TextValue $ws.Cells.Item(5,4) '558.80'
$ws.Cells.Item(5,5).Value2 = '  -0.49%  '

Set-TextValue $ws.Cells.Item(6,4) '160.86'
$ws.Cells.Item(6,5).Value2 = '  -0.92%  '

$ws.Cells.Item(7,5).Value2 = '  +0.00%  '

Set-TextValue $ws.Cells.Item(8,4) '0.510'
$ws.Cells.Item(8,5).Value2 = '  +0.76%  '

Set-TextValue $ws.Cells.Item(9,4) '0.167'
$ws.Cells.Item(9,5).Value2 = '  +10.08%  '

$ws.Cells.Item(10,5).Value2 = '  -1.69%  '

Set-TextValue $ws.Cells.Item(11,4) '0.330'
$ws.Cells.Item(11,5).Value2 = '  -0.83%  '

$ws.Cells.Item(12,5).Value2 = '  -5.68%  '

$ws.Cells.Item(13,2).Value2 = 'WrappedBTC'
$ws.Cells.Item(13,3).Value2 = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue $ws.Cells.Item(13,4) '68.628.06'
$ws.Cells.Item(13,5).Value2 = '  -0.42%  '

$ws.Cells.Item(14,2).Value2 = 'ShibaInu'
$ws.Cells.Item(14,3).Value2 = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextValue $ws.Cells.Item(14,4) '0.0000175'
$ws.Cells.Item(14,5).Value2 = '  +3.67%  '

Set-TextValue $ws.Cells.Item(15,4) '2.872.56'
$ws.Cells.Item(15,5).Value2 = '  -1.14%  '

$ws.Cells.Item(16,5).Value2 = '  -2.23%  '

Set-TextValue $ws.Cells.Item(17,4) '2.427.54'
$ws.Cells.Item(17,5).Value2 = '  -1.43%  '

$ws.Cells.Item(18,5).Value2 = '  -2.32%  '

Set-TextValue $ws.Cells.Item(19,4) '334.58'
$ws.Cells.Item(19,5).Value2 = '  -0.66%  '

Set-TextValue $ws.Cells.Item(20,4) '6.91'
$ws.Cells.Item(20,5).Value2 = '  -0.87%  '

$ws.Cells.Item(21,5).Value2 = '  +0.27%  '

$ws.Cells.Item(22,5).Value2 = '  +2.66%  '

Set-TextValue $ws.Cells.Item(23,4) '1.00'
$ws.Cells.Item(23,5).Value2 = '  +0.17%  '

Set-TextValue $ws.Cells.Item(24,4) '66.92'
$ws.Cells.Item(24,5).Value2 = '  +0.11%  '

Set-TextValue $ws.Cells.Item(25,4) '3.68'
$ws.Cells.Item(25,5).Value2 = '  +0.12%  '

Set-TextValue $ws.Cells.Item(26,4) '2.554.46'

Set-TextValue $ws.Cells.Item(27,4) '1.01'
$ws.Cells.Item(27,5).Value2 = '  +0.86%  '

Set-TextValue $ws.Cells.Item(28,4) '8.18'
$ws.Cells.Item(28,5).Value2 = '  -0.56%  '

$ws.Cells.Item(29,5).Value2 = '  -0.54%  '

$ws.Cells.Item(30,5).Value2 = '  -1.59%  '

$ws.Cells.Item(31,5).Value2 = '  +0.11%  '

Set-TextValue $ws.Cells.Item(32,4) '427.49'
$ws.Cells.Item(32,5).Value2 = '  -0.79%  '

$ws.Cells.Item(33,5).Value2 = '  +0.85%  '

$ws.Cells.Item(34,5).Value2 = '  -0.49%  '

Set-TextValue $ws.Cells.Item(35,4) '160.65'
$ws.Cells.Item(35,5).Value2 = '  +1.44%  '

$ws.Cells.Item(37,5).Value2 = '  +0.00%  '

Set-TextValue $ws.Cells.Item(38,4) '17.91'
$ws.Cells.Item(38,5).Value2 = '  +0.60%  '

$ws.Cells.Item(40,5).Value2 = '  -1.26%  '

$ws.Cells.Item(41,2).Value2 = 'Stacks'
$ws.Cells.Item(41,3).Value2 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue $ws.Cells.Item(41,4) '1.50'
$ws.Cells.Item(41,5).Value2 = '  +1.61%  '

$ws.Cells.Item(42,2).Value2 = 'RenderToken'
$ws.Cells.Item(42,3).Value2 = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue $ws.Cells.Item(42,4) '4.34'
$ws.Cells.Item(42,5).Value2 = '  -2.15%  '

$ws.Cells.Item(43,5).Value2 = '  -0.74%  '

$ws.Cells.Item(44,5).Value2 = '  -1.81%  '

Set-TextValue $ws.Cells.Item(45,4) '131.91'
$ws.Cells.Item(45,5).Value2 = '  +0.17%  '

$ws.Cells.Item(46,5).Value2 = '  -0.58%  '

Set-TextValue $ws.Cells.Item(47,4) '0.0716'
$ws.Cells.Item(47,5).Value2 = '  +0.47%  '

$ws.Cells.Item(48,5).Value2 = '  -0.98%  '

$ws.Cells.Item(49,5).Value2 = '  -1.22%  '

$ws.Cells.Item(50,5).Value2 = '  +0.51%  '

$ws.Cells.Item(51,5).Value2 = '  +0.44%  '
